$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G62").Copy()
$ws.Range("G14").PasteSpecial(-4104)
